# Update Leave Card 8/3/2023 4:40 PM
# Records two new SL(1-0-0) sick-leave entries (7/21/2023 and 7/31/2023) in
# the leave card table, adds the earned credit for the existing SP(1-0-0)
# entry, and grows the table by the two extra rows this creates (shifting
# the recurring PERIOD dates and the trailing template rows down by 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table by 2 rows up front so later writes/copies into the
#        newly-appended rows are recognised as part of Table1. ---
$lo.Resize($ws.Range("A8:K136"))

# --- 2. Push the two trailing "end of table" rows (133,134) down to
#        (135,136), then rebuild 133/134 as normal data rows (same layout
#        as row 132). ---
$ws.Range("A133:K134").Copy($ws.Range("A135:K136"))
$ws.Range("A132:K132").Copy($ws.Range("A133:K133"))
$ws.Range("A132:K132").Copy($ws.Range("A134:K134"))

# Restore the literal calculated-column formula text (Copy() round-trips it
# fine for 133/134, but keep 135/136 explicit too for safety/clarity).
$earnedFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G133").Formula = $earnedFormula
$ws.Range("G134").Formula = $earnedFormula
$ws.Range("G135").Formula = $earnedFormula
$ws.Range("G136").Formula = $earnedFormula

# --- 3. Make room for the two new leave-card rows at 27/28: shift every
#        recurring PERIOD date in column A down by 2 rows (old row r -> new
#        row r+2), from the bottom up so destinations are empty first. ---
for ($r = 121; $r -ge 27; $r--) {
    $srcVal = $ws.Range("A$r").Value2
    $dest = $r + 2
    $ws.Range("A$dest").Value = $srcVal
}
$ws.Range("A27").ClearContents()
$ws.Range("A28").ClearContents()

# --- 4. Fill in the EARNED credit for the existing SP(1-0-0) row (26). ---
$ws.Range("C26").Value = 1.25

# --- 5. New SL(1-0-0) entry used on 7/21/2023 (row 27). ---
$ws.Range("B27").Value = "SL(1-0-0)"
$ws.Range("H27").Value = 1
$ws.Range("K26").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 45128

# --- 6. New SL(1-0-0) entry used on 7/31/2023 (row 28). ---
$ws.Range("B28").Value = "SL(1-0-0)"
$ws.Range("H28").Value = 1
$ws.Range("K26").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = 45138

# --- 7. Move the active selection to where the user last clicked. ---
$ws.Range("I28").Select()
